# Update the "participants" sheet field list:
#  - remove the "date" and "age_days" columns
#  - remove the "consent" column
#  - add a new "birth weight" column at the end
#  - reorder the remaining fields
#  - fix up the data validation ranges to match the new column layout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("participants")

# Drop every existing data validation rule up front; they get re-created below
# against the final column layout, so there's no need to track how each one
# drifts as columns are removed.
$ws.Range("A1:Z1001").Validation.Delete()

# Drop the now-unused "date" column (old column C) and "age_days" column
# (old column D, which shifts into C once the first delete happens).
$ws.Columns("C").Delete()
$ws.Columns("C").Delete()

# Drop the "consent" column (old column K, now column I after the two deletes above).
$ws.Columns("I").Delete()

# Rewrite the header row to match the new target order and append "birth weight".
$ws.Range("A1").Value = "participantID"
$ws.Range("B1").Value = "birthdate"
$ws.Range("C1").Value = "gender"
$ws.Range("D1").Value = "race"
$ws.Range("E1").Value = "ethnicity"
$ws.Range("F1").Value = "language"
$ws.Range("G1").Value = "disability"
$ws.Range("H1").Value = "gestational age"
$ws.Range("I1").Value = "pregnancy term"
$ws.Range("J1").Value = "birth weight"

# Re-create the data validation rules against the new columns.
$ws.Range("C2:C1001").Validation.Add(3, 1, 1, '"Female,Male"')
$ws.Range("D2:D1001").Validation.Add(3, 1, 1, '"American Indian or Alaska Native,Asian,Native Hawaiian or Other Pacific Islander,Black or African American,White,More than one,Unknown or not reported"')
$ws.Range("E2:E1001").Validation.Add(3, 1, 1, '"Not Hispanic or Latino,Hispanic or Latino"')
$ws.Range("I2:I1001").Validation.Add(3, 1, 1, '"Full term,Preterm"')

Write-Host "done"
